$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 429, shifting existing rows 429:505 down to 430:506
$ws.Rows.Item(429).Insert()

# Populate the newly inserted row 429 with the new record
$ws.Cells.Item(429, 1).Value = 10
$ws.Cells.Item(429, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(429, 3).Value = "La Araucanía"
$ws.Cells.Item(429, 4).Value = 44889
$ws.Cells.Item(429, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(429, 5).Value = 9
$ws.Cells.Item(429, 6).Value = 100112008
$ws.Cells.Item(429, 7).Value = "Coliflor"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 1250
$ws.Cells.Item(429, 11).Value = 1100
$ws.Cells.Item(429, 12).Value = 1100
$ws.Cells.Item(429, 13).Value = 1100
$ws.Cells.Item(429, 14).Value = "$/unidad"
$ws.Cells.Item(429, 15).Value = "Región del Maule"
$ws.Cells.Item(429, 16).Value = 1100
$ws.Cells.Item(429, 17).Value = 1
$ws.Cells.Item(429, 18).Value = "Hortaliza"
